# Auto-generated edit script: apply updated market-price / profit values
# to the Atomos_Profits leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 886.55554
$ws.Range("I6").Value = 163.16667
$ws.Range("J6").Value = 2333.3333
$ws.Range("K6").Value = 489.50001
$ws.Range("L6").Value = 6999.999899999999
$ws.Range("M6").Value = -377.50001
$ws.Range("N6").Value = -7223.999899999999
$ws.Range("H74").Value = 4213
$ws.Range("I74").Value = 7001.5
$ws.Range("J74").Value = 3593.3333
$ws.Range("K74").Value = 7001.5
$ws.Range("L74").Value = 3593.3333
$ws.Range("M74").Value = -6065.5
$ws.Range("N74").Value = -5465.3333
$ws.Range("H77").Value = 4213
$ws.Range("I77").Value = 7001.5
$ws.Range("J77").Value = 3593.3333
$ws.Range("K77").Value = 35007.5
$ws.Range("L77").Value = 17966.6665
$ws.Range("M77").Value = -30327.5
$ws.Range("N77").Value = -27326.6665
$ws.Range("H115").Value = 1604.4166
$ws.Range("I115").Value = 650.6
$ws.Range("J115").Value = 2285.7144
$ws.Range("K115").Value = 1951.8
$ws.Range("L115").Value = 6857.1432
$ws.Range("M115").Value = -384.8000000000002
$ws.Range("N115").Value = -9991.143199999999
$ws.Range("H129").Value = 6098820
$ws.Range("J129").Value = 1301.9395
$ws.Range("L129").Value = 3905.8185
$ws.Range("N129").Value = -13905.8185
$ws.Range("H134").Value = 32000
$ws.Range("J134").Value = 32000
$ws.Range("L134").Value = 32000
$ws.Range("N134").Value = -42140
$ws.Range("H140").Value = 48312
$ws.Range("J140").Value = 48312
$ws.Range("L140").Value = 48312
$ws.Range("N140").Value = -58672

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 40712
$ws.Range("J7").Value = 40712
$ws.Range("L7").Value = 40712
$ws.Range("N7").Value = -40940
$ws.Range("H52").Value = 39780
$ws.Range("J52").Value = 39780
$ws.Range("L52").Value = 39780
$ws.Range("N52").Value = -40416
$ws.Range("H62").Value = 4500
$ws.Range("J62").Value = 4500
$ws.Range("L62").Value = 4500
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 4500
$ws.Range("J65").Value = 4500
$ws.Range("L65").Value = 13500
$ws.Range("N65").Value = -19740
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H130").Value = 36000
$ws.Range("J130").Value = 36000
$ws.Range("L130").Value = 36000
$ws.Range("N130").Value = -46040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 30048.75
$ws.Range("J45").Value = 30048.75
$ws.Range("L45").Value = 30048.75
$ws.Range("N45").Value = -31664.75
$ws.Range("H86").Value = 618424.8
$ws.Range("I86").Value = 1150998.8
$ws.Range("J86").Value = 3916.3845
$ws.Range("K86").Value = 1150998.8
$ws.Range("L86").Value = 3916.3845
$ws.Range("M86").Value = -1149875.8
$ws.Range("N86").Value = -6162.3845
$ws.Range("H89").Value = 618424.8
$ws.Range("I89").Value = 1150998.8
$ws.Range("J89").Value = 3916.3845
$ws.Range("K89").Value = 5754994
$ws.Range("L89").Value = 19581.9225
$ws.Range("M89").Value = -5749378
$ws.Range("N89").Value = -30813.9225
$ws.Range("H108").Value = 36000
$ws.Range("J108").Value = 36000
$ws.Range("L108").Value = 36000
$ws.Range("N108").Value = -43680
$ws.Range("H134").Value = 3580.1177
$ws.Range("I134").Value = 2571.8333
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 7715.499899999999
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -5180.499899999999
$ws.Range("N134").Value = -23070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9436191
$ws.Range("I58").Value = 939.1111
$ws.Range("K58").Value = 939.1111
$ws.Range("M58").Value = -736.1111
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
$ws.Range("H120").Value = 20100
$ws.Range("J120").Value = 20100
$ws.Range("L120").Value = 20100
$ws.Range("N120").Value = -27358
$ws.Range("H132").Value = 2106.889
$ws.Range("I132").Value = 1599
$ws.Range("J132").Value = 3261.182
$ws.Range("K132").Value = 4797
$ws.Range("L132").Value = 9783.545999999998
$ws.Range("M132").Value = -2267
$ws.Range("N132").Value = -14843.546
$ws.Range("H136").Value = 9436191
$ws.Range("I136").Value = 939.1111
$ws.Range("K136").Value = 2817.3333
$ws.Range("M136").Value = -267.3332999999998
$ws.Range("H141").Value = 24760.857
$ws.Range("J141").Value = 24760.857
$ws.Range("L141").Value = 24760.857
$ws.Range("N141").Value = -35120.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2795.2
$ws.Range("J55").Value = 3237.1428
$ws.Range("L55").Value = 9711.428400000001
$ws.Range("N55").Value = -10065.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2321.2173
$ws.Range("I97").Value = 1561.2354
$ws.Range("J97").Value = 4474.5
$ws.Range("K97").Value = 1561.2354
$ws.Range("L97").Value = 4474.5
$ws.Range("M97").Value = -1065.2354
$ws.Range("N97").Value = -5466.5
$ws.Range("H126").Value = 2932.68
$ws.Range("I126").Value = 1685.9231
$ws.Range("J126").Value = 4283.3335
$ws.Range("K126").Value = 5057.7693
$ws.Range("L126").Value = 12850.0005
$ws.Range("M126").Value = -2587.7693
$ws.Range("N126").Value = -17790.0005
$ws.Range("H129").Value = 41249.625
$ws.Range("J129").Value = 41249.625
$ws.Range("L129").Value = 41249.625
$ws.Range("N129").Value = -51249.625
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
$ws.Range("H135").Value = 28833.334
$ws.Range("J135").Value = 28833.334
$ws.Range("L135").Value = 28833.334
$ws.Range("N135").Value = -38973.334
$ws.Range("H141").Value = 46833
$ws.Range("J141").Value = 46833
$ws.Range("L141").Value = 46833
$ws.Range("N141").Value = -57193

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 90909990
$ws.Range("I16").Value = 125001144
$ws.Range("J16").Value = 248
$ws.Range("K16").Value = 125001144
$ws.Range("L16").Value = 248
$ws.Range("M16").Value = -125000974
$ws.Range("N16").Value = -588
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 22900
$ws.Range("J131").Value = 22900
$ws.Range("L131").Value = 22900
$ws.Range("N131").Value = -32980
$ws.Range("H132").Value = 1559.871
$ws.Range("I132").Value = 884.2041
$ws.Range("J132").Value = 4106.615
$ws.Range("K132").Value = 2652.6123
$ws.Range("L132").Value = 12319.845
$ws.Range("M132").Value = -122.6123000000002
$ws.Range("N132").Value = -17379.845
$ws.Range("H135").Value = 29551.25
$ws.Range("J135").Value = 29551.25
$ws.Range("L135").Value = 29551.25
$ws.Range("N135").Value = -39691.25
$ws.Range("H136").Value = 1691.5588
$ws.Range("I136").Value = 1072.1072
$ws.Range("J136").Value = 4582.3335
$ws.Range("K136").Value = 3216.3216
$ws.Range("L136").Value = 13747.0005
$ws.Range("M136").Value = -666.3215999999998
$ws.Range("N136").Value = -18847.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1147.45
$ws.Range("I100").Value = 1100.625
$ws.Range("K100").Value = 2201.25
$ws.Range("M100").Value = -1660.25
$ws.Range("H113").Value = 2135.5454
$ws.Range("I113").Value = 618.4
$ws.Range("K113").Value = 1855.2
$ws.Range("M113").Value = 314.8000000000002
$ws.Range("H132").Value = 26503.695
$ws.Range("I132").Value = 5036.375
$ws.Range("J132").Value = 75571.86
$ws.Range("K132").Value = 15109.125
$ws.Range("L132").Value = 226715.58
$ws.Range("M132").Value = -12579.125
$ws.Range("N132").Value = -231775.58
$ws.Range("H137").Value = 46357.5
$ws.Range("J137").Value = 46357.5
$ws.Range("L137").Value = 46357.5
$ws.Range("N137").Value = -56557.5
$ws.Range("H140").Value = 63332.332
$ws.Range("J140").Value = 63332.332
$ws.Range("L140").Value = 63332.332
$ws.Range("N140").Value = -73692.33199999999
$ws.Range("H141").Value = 28388.889
$ws.Range("J141").Value = 28388.889
$ws.Range("L141").Value = 28388.889
$ws.Range("N141").Value = -38748.889
